$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column L = "Hintergrundmodell" flag, values "ja" / "nein"
$ws.Range("L4").Value = "ja"
$ws.Range("L9").Value = "ja"
$ws.Range("L12").Value = "nein"
$ws.Range("L13").Value = "nein"
$ws.Range("L14").Value = "nein"
$ws.Range("L15").Value = "nein"
